$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.366.13"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.950.31"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.79%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.42"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.612"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.71"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -7.40%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.369"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -5.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "55.64"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0822"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.02%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.828"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -8.19%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.63"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -6.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.238.38"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.48"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -5.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.30"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.938.93"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.236.53"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.75"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0869"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.26"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.01"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -6.67%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.50"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.27"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.46"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.21"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.53"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.118"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -14.98%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.06%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.72"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -6.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0634"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.30"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.38%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.07"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.80"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.14"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -9.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.88"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -8.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0981"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.88"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.18"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0209"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.63"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -8.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.03"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -8.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.41"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.342.20"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.34"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.58%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.91"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.40%  "
